$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.218295574188232
$ws.Range("B1").Value = 2.832104921340942
$ws.Range("C1").Value = 2.176068544387817
$ws.Range("D1").Value = 2.033602237701416
$ws.Range("E1").Value = 2.04894495010376
